$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 09:04"

# Row 38: Ucrania - refreshed case counts
$ws.Range("B38").Value = 13691
$ws.Range("C38").Value = 507
$ws.Range("D38").Value = 2396
$ws.Range("E38").Value = 10955
$ws.Range("F38").Value = 178
$ws.Range("G38").Value = 13
$ws.Range("H38").Value = 340

# Row 94: Letonia - refreshed case counts
$ws.Range("B94").Value = 909
$ws.Range("C94").Value = 9
$ws.Range("E94").Value = 427
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 18

# Row 98: was Sudan, now Consejo Danes para los Refugiados (new data causes it to rank here)
$ws.Range("A98").Value = "Consejo Danes para los Refugiados"
$ws.Range("B98").Value = 863
$ws.Range("C98").Value = 66
$ws.Range("D98").Value = 103
$ws.Range("E98").Value = 724
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 36

# Row 99: was Albania, now Sudan (shifted down one rank)
$ws.Range("A99").Value = "Sudan"
$ws.Range("B99").Value = 852
$ws.Range("D99").Value = 80
$ws.Range("E99").Value = 723
$ws.Range("F99").Value = 0
$ws.Range("H99").Value = 49

# Row 100: was Guatemala, now Albania (shifted down one rank)
$ws.Range("A100").Value = "Albania"
$ws.Range("B100").Value = 832
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 595
$ws.Range("E100").Value = 206
$ws.Range("F100").Value = 7
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 31

# Row 101: was Consejo Danes para los Refugiados, now Guatemala (shifted down one rank)
$ws.Range("A101").Value = "Guatemala"
$ws.Range("B101").Value = 798
$ws.Range("C101").Value = 35
$ws.Range("D101").Value = 86
$ws.Range("E101").Value = 691
$ws.Range("F101").Value = 5
$ws.Range("G101").Value = 2
$ws.Range("H101").Value = 21

# Row 191: was Belice, now Nueva Caledonia (swapped rank with row 192)
$ws.Range("A191").Value = "Nueva Caledonia"
$ws.Range("D191").Value = 18
$ws.Range("H191").Value = 0

# Row 192: was Nueva Caledonia, now Belice (swapped rank with row 191)
$ws.Range("A192").Value = "Belice"
$ws.Range("D192").Value = 16
$ws.Range("H192").Value = 2

# Row 198: was Curazao, now Dominica (swapped rank with row 199)
$ws.Range("A198").Value = "Dominica"
$ws.Range("D198").Value = 14
$ws.Range("H198").Value = 0

# Row 199: was Dominica, now Curazao (swapped rank with row 198)
$ws.Range("A199").Value = "Curazao"
$ws.Range("D199").Value = 13
$ws.Range("H199").Value = 1
